# Internationalize the worksheet contents: translate Chinese header/labels
# and name/department values into their English equivalents.
# Only the text (string) cell values change; numeric cell values are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Department"
$ws.Range("D1").Value = "Salary"
$ws.Range("E1").Value = "Performance Score"

# Row 2 - Zhang San
$ws.Range("A2").Value = "Zhang San"
$ws.Range("C2").Value = "Technology"

# Row 3 - Li Si
$ws.Range("A3").Value = "Li Si"
$ws.Range("C3").Value = "Sales"

# Row 4 - Wang Wu
$ws.Range("A4").Value = "Wang Wu"
$ws.Range("C4").Value = "Marketing"

# Row 5 - Zhao Liu
$ws.Range("A5").Value = "Zhao Liu"
$ws.Range("C5").Value = "Technology"

# Row 6 - Qian Qi
$ws.Range("A6").Value = "Qian Qi"
$ws.Range("C6").Value = "HR"

# Row 7 - Sun Ba
$ws.Range("A7").Value = "Sun Ba"
$ws.Range("C7").Value = "Sales"

# Row 8 - Zhou Jiu
$ws.Range("A8").Value = "Zhou Jiu"
$ws.Range("C8").Value = "Technology"

# Row 9 - Wu Shi
$ws.Range("A9").Value = "Wu Shi"
$ws.Range("C9").Value = "Marketing"

# Row 10 - Zheng Shiyi
$ws.Range("A10").Value = "Zheng Shiyi"
$ws.Range("C10").Value = "Technology"

# Row 11 - Wang Shier
$ws.Range("A11").Value = "Wang Shier"
$ws.Range("C11").Value = "Sales"
